# Updates cryptocurrency price/volume/hour data to match the
# Jan 6 2023 21:02 UTC GitHub Actions refresh of the symbol list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the target cells as Text ("@") so values like "258.97",
# "0.45%" and "21" are stored verbatim instead of being coerced into
# numbers / percentages by Excel's normal type inference. Each
# contiguous block is set individually -- multi-area (comma) ranges
# only apply NumberFormat to their first area in this host.
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D8").NumberFormat = "@"
$ws.Range("D10:D11").NumberFormat = "@"
$ws.Range("D13:D18").NumberFormat = "@"
$ws.Range("D20:D24").NumberFormat = "@"
$ws.Range("D26:D27").NumberFormat = "@"
$ws.Range("D40:D45").NumberFormat = "@"
$ws.Range("D47:D50").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E4:E27").NumberFormat = "@"
$ws.Range("E40:E50").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Column D (Price), Column E (Volume(1h)), Column G (Hora) updates
$ws.Range("D2").Value = '258.97'
$ws.Range("E2").Value = '0.45%'
$ws.Range("G2").Value = '21'
$ws.Range("D3").Value = '26.95'
$ws.Range("G3").Value = '21'
$ws.Range("E4").Value = '0.47%'
$ws.Range("G4").Value = '21'
$ws.Range("D5").Value = '0.06033'
$ws.Range("E5").Value = '2.49%'
$ws.Range("G5").Value = '21'
$ws.Range("D6").Value = '6.681'
$ws.Range("E6").Value = '0.52%'
$ws.Range("G6").Value = '21'
$ws.Range("D7").Value = '0.8576'
$ws.Range("E7").Value = '-0.02%'
$ws.Range("G7").Value = '21'
$ws.Range("D8").Value = '0.9240'
$ws.Range("E8").Value = '-1.99%'
$ws.Range("G8").Value = '21'
$ws.Range("E9").Value = '-0.79%'
$ws.Range("G9").Value = '21'
$ws.Range("D10").Value = '0.05016'
$ws.Range("E10").Value = '30.50%'
$ws.Range("G10").Value = '21'
$ws.Range("D11").Value = '0.07035'
$ws.Range("E11").Value = '-0.73%'
$ws.Range("G11").Value = '21'
$ws.Range("E12").Value = '-1.59%'
$ws.Range("G12").Value = '21'
$ws.Range("D13").Value = '0.09132'
$ws.Range("E13").Value = '-0.56%'
$ws.Range("G13").Value = '21'
$ws.Range("D14").Value = '0.001541'
$ws.Range("E14").Value = '-0.14%'
$ws.Range("G14").Value = '21'
$ws.Range("D15").Value = '0.0006078'
$ws.Range("E15").Value = '0.63%'
$ws.Range("G15").Value = '21'
$ws.Range("D16").Value = '0.006036'
$ws.Range("E16").Value = '-2.97%'
$ws.Range("G16").Value = '21'
$ws.Range("D17").Value = '3.462'
$ws.Range("E17").Value = '-1.47%'
$ws.Range("G17").Value = '21'
$ws.Range("D18").Value = '3.171'
$ws.Range("E18").Value = '-1.09%'
$ws.Range("G18").Value = '21'
$ws.Range("E19").Value = '-1.74%'
$ws.Range("G19").Value = '21'
$ws.Range("D20").Value = '0.3093'
$ws.Range("E20").Value = '0.42%'
$ws.Range("G20").Value = '21'
$ws.Range("D21").Value = '0.1298'
$ws.Range("E21").Value = '0.45%'
$ws.Range("G21").Value = '21'
$ws.Range("D22").Value = '4.128'
$ws.Range("E22").Value = '6.58%'
$ws.Range("G22").Value = '21'
$ws.Range("D23").Value = '0.04234'
$ws.Range("E23").Value = '0.19%'
$ws.Range("G23").Value = '21'
$ws.Range("D24").Value = '0.001218'
$ws.Range("E24").Value = '-0.21%'
$ws.Range("G24").Value = '21'
$ws.Range("E25").Value = '-6.06%'
$ws.Range("G25").Value = '21'
$ws.Range("D26").Value = '0.0001201'
$ws.Range("E26").Value = '0.10%'
$ws.Range("G26").Value = '21'
$ws.Range("D27").Value = '0.0001524'
$ws.Range("E27").Value = '-21.29%'
$ws.Range("G27").Value = '21'
$ws.Range("G28").Value = '21'
$ws.Range("G29").Value = '21'
$ws.Range("G30").Value = '21'
$ws.Range("G31").Value = '21'
$ws.Range("G32").Value = '21'
$ws.Range("G33").Value = '21'
$ws.Range("G34").Value = '21'
$ws.Range("G35").Value = '21'
$ws.Range("G36").Value = '21'
$ws.Range("G37").Value = '21'
$ws.Range("G38").Value = '21'
$ws.Range("G39").Value = '21'
$ws.Range("D40").Value = '0.03842'
$ws.Range("E40").Value = '0.33%'
$ws.Range("G40").Value = '21'
$ws.Range("D41").Value = '0.1115'
$ws.Range("E41").Value = '1.14%'
$ws.Range("G41").Value = '21'
$ws.Range("D42").Value = '0.003956'
$ws.Range("E42").Value = '-36.17%'
$ws.Range("G42").Value = '21'
$ws.Range("D43").Value = '0.01525'
$ws.Range("E43").Value = '33.26%'
$ws.Range("G43").Value = '21'
$ws.Range("D44").Value = '0.002201'
$ws.Range("E44").Value = '0.10%'
$ws.Range("G44").Value = '21'
$ws.Range("D45").Value = '0.00005083'
$ws.Range("E45").Value = '-6.86%'
$ws.Range("G45").Value = '21'
$ws.Range("E46").Value = '0.10%'
$ws.Range("G46").Value = '21'
$ws.Range("D47").Value = '0.05458'
$ws.Range("E47").Value = '-9.02%'
$ws.Range("G47").Value = '21'
$ws.Range("D48").Value = '0.1321'
$ws.Range("E48").Value = '2.30%'
$ws.Range("G48").Value = '21'
$ws.Range("D49").Value = '0.00002101'
$ws.Range("E49").Value = '0.10%'
$ws.Range("G49").Value = '21'
$ws.Range("D50").Value = '0.0002001'
$ws.Range("E50").Value = '0.10%'
$ws.Range("G50").Value = '21'
$ws.Range("G51").Value = '21'

Write-Output "Updated 117 cells across D/E/G columns"
